$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("M1").Value = "Future Dream"
$ws.Range("N1").Value = "Difficulty"

# Copy header style (bold + border) from existing header cell L1 to M1:N1
$ws.Range("L1").Copy() | Out-Null
$ws.Range("M1:N1").PasteSpecial(-4122) | Out-Null

$data = @(
    @("[0, 1, 0, 0, 0, 0]", "[0, 1, 0]"),
    @("[0, 0, 1, 0, 0, 0]", "[1, 0, 0]"),
    @("[1, 0, 0, 0, 0, 0]", "[0, 1, 0]"),
    @("[0, 0, 0, 0, 1, 0]", "[1, 0, 0]"),
    @("[0, 0, 0, 1, 0, 0]", "[1, 0, 0]"),
    @("[0, 1, 0, 0, 0, 0]", "[0, 1, 0]"),
    @("[0, 0, 1, 0, 0, 0]", "[1, 0, 0]"),
    @("[0, 0, 0, 0, 0, 1]", "[0, 0, 1]")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 13).Value = $data[$i][0]
    $ws.Cells.Item($row, 14).Value = $data[$i][1]
}
